$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 889011.9
$ws.Range("I19").Value = 1212205.4
$ws.Range("K19").Value = 1212205.4
$ws.Range("M19").Value = -1212030.4

$ws.Range("H64").Value = 2609
$ws.Range("I64").Value = 2550
$ws.Range("J64").Value = 2679.8
$ws.Range("K64").Value = 2550
$ws.Range("L64").Value = 2679.8
$ws.Range("M64").Value = -2302
$ws.Range("N64").Value = -3175.8

$ws.Range("H67").Value = 2609
$ws.Range("I67").Value = 2550
$ws.Range("J67").Value = 2679.8
$ws.Range("K67").Value = 2550
$ws.Range("L67").Value = 2679.8
$ws.Range("M67").Value = -1692
$ws.Range("N67").Value = -4395.8

$ws.Range("H76").Value = 3246.6
$ws.Range("I76").Value = 3235.6428
$ws.Range("J76").Value = 3400
$ws.Range("K76").Value = 3235.6428
$ws.Range("L76").Value = 3400
$ws.Range("M76").Value = -2920.6428
$ws.Range("N76").Value = -4030

$ws.Range("H79").Value = 3246.6
$ws.Range("I79").Value = 3235.6428
$ws.Range("J79").Value = 3400
$ws.Range("K79").Value = 3235.6428
$ws.Range("L79").Value = 3400
$ws.Range("M79").Value = -2143.6428
$ws.Range("N79").Value = -5584

$ws.Range("H112").Value = 10102575
$ws.Range("J112").Value = 1597.125
$ws.Range("L112").Value = 4791.375
$ws.Range("N112").Value = -7007.375

$ws.Range("H129").Value = 30612.906
$ws.Range("J129").Value = 32627.8
$ws.Range("L129").Value = 97883.39999999999
$ws.Range("N129").Value = -107883.4

$ws.Range("H132").Value = 91967.84
$ws.Range("I132").Value = 116273.23
$ws.Range("K132").Value = 348819.69
$ws.Range("M132").Value = -346289.69

$ws.Range("H135").Value = 1039.3334
$ws.Range("I135").Value = 771.0909
$ws.Range("J135").Value = 3990
$ws.Range("K135").Value = 6939.8181
$ws.Range("L135").Value = 35910
$ws.Range("M135").Value = -4404.8181
$ws.Range("N135").Value = -40980

$ws.Range("H137").Value = 6141.64
$ws.Range("I137").Value = 6252.533
$ws.Range("K137").Value = 18757.599
$ws.Range("M137").Value = -16207.599

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3276.4
$ws.Range("I74").Value = 3311.2778
$ws.Range("K74").Value = 3311.2778
$ws.Range("M74").Value = -2437.2778

$ws.Range("H77").Value = 3276.4
$ws.Range("I77").Value = 3311.2778
$ws.Range("K77").Value = 16556.389
$ws.Range("M77").Value = -12188.389

$ws.Range("H132").Value = 2859.7317
$ws.Range("I132").Value = 2079.889
$ws.Range("J132").Value = 4363.7144
$ws.Range("K132").Value = 6239.667
$ws.Range("L132").Value = 13091.1432
$ws.Range("M132").Value = -3709.667
$ws.Range("N132").Value = -18151.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2542.8572
$ws.Range("I86").Value = 2560
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 2560
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -1437
$ws.Range("N86").Value = -4746

$ws.Range("H89").Value = 2542.8572
$ws.Range("I89").Value = 2560
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 12800
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -7184
$ws.Range("N89").Value = -23732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 51.666668
$ws.Range("I7").Value = 40.833332
$ws.Range("J7").Value = 73.333336
$ws.Range("K7").Value = 40.833332
$ws.Range("L7").Value = 73.333336
$ws.Range("M7").Value = 72.166668
$ws.Range("N7").Value = -299.333336

$ws.Range("H58").Value = 2165.1228
$ws.Range("I58").Value = 1807.9231
$ws.Range("K58").Value = 1807.9231
$ws.Range("M58").Value = -1604.9231

$ws.Range("H62").Value = 62503310
$ws.Range("I62").Value = 62503310
$ws.Range("K62").Value = 62503310
$ws.Range("M62").Value = -62502686

$ws.Range("H65").Value = 62503310
$ws.Range("I65").Value = 62503310
$ws.Range("K65").Value = 312516550
$ws.Range("M65").Value = -312513430

$ws.Range("H136").Value = 2165.1228
$ws.Range("I136").Value = 1807.9231
$ws.Range("K136").Value = 5423.7693
$ws.Range("M136").Value = -2873.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 27000
$ws.Range("I58").Value = 27000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 27000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -26723
$ws.Range("N58").ClearContents()

$ws.Range("H70").Value = 6430.853
$ws.Range("I70").Value = 5763.423
$ws.Range("J70").Value = 8600
$ws.Range("K70").Value = 5763.423
$ws.Range("L70").Value = 8600
$ws.Range("M70").Value = -5493.423
$ws.Range("N70").Value = -9140

$ws.Range("H73").Value = 6430.853
$ws.Range("I73").Value = 5763.423
$ws.Range("J73").Value = 8600
$ws.Range("K73").Value = 5763.423
$ws.Range("L73").Value = 8600
$ws.Range("M73").Value = -4827.423
$ws.Range("N73").Value = -10472

$ws.Range("H102").Value = 3008.923
$ws.Range("I102").Value = 2335.6667
$ws.Range("J102").Value = 3927
$ws.Range("K102").Value = 2335.6667
$ws.Range("L102").Value = 3927
$ws.Range("M102").Value = -713.6667000000002
$ws.Range("N102").Value = -7171

$ws.Range("H132").Value = 2901.24
$ws.Range("I132").Value = 1159.5385
$ws.Range("J132").Value = 4788.0835
$ws.Range("K132").Value = 3478.6155
$ws.Range("L132").Value = 14364.2505
$ws.Range("M132").Value = -948.6155000000003
$ws.Range("N132").Value = -19424.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7036.625
$ws.Range("I132").Value = 2480.8
$ws.Range("J132").Value = 8235.526
$ws.Range("K132").Value = 7442.400000000001
$ws.Range("L132").Value = 24706.578
$ws.Range("M132").Value = -4912.400000000001
$ws.Range("N132").Value = -29766.578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 23810856
$ws.Range("I81").Value = 23810856
$ws.Range("K81").Value = 47621712
$ws.Range("M81").Value = -47620651

$ws.Range("H84").Value = 23810856
$ws.Range("I84").Value = 23810856
$ws.Range("K84").Value = 238108560
$ws.Range("M84").Value = -238103256

$ws.Range("H122").Value = 4639.75
$ws.Range("I122").Value = 3497
$ws.Range("J122").Value = 5574.727
$ws.Range("K122").Value = 10491
$ws.Range("L122").Value = 16724.181
$ws.Range("M122").Value = -8041
$ws.Range("N122").Value = -21624.181

$ws.Range("H132").Value = 9261296
$ws.Range("I132").Value = 781.38464
$ws.Range("J132").Value = 14495500
$ws.Range("K132").Value = 2344.15392
$ws.Range("L132").Value = 43486500
$ws.Range("M132").Value = 185.8460800000003
$ws.Range("N132").Value = -43491560

$ws.Range("H136").Value = 2739.6191
$ws.Range("I136").Value = 1311
$ws.Range("J136").Value = 5596.857
$ws.Range("K136").Value = 3933
$ws.Range("L136").Value = 16790.571
$ws.Range("M136").Value = -1383
$ws.Range("N136").Value = -21890.571
